$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new student record in row 11 (practice properties file read/write data)
$ws.Range("A11").Value = 101
$ws.Range("B11").Value = "dfa"
$ws.Range("C11").Value = 33
$ws.Range("D11").Value = 32
$ws.Range("E11").Value = 32
$ws.Range("F11").Value = 44
$ws.Range("G11").Value = "C"
$ws.Range("H11").Value = "fasd"
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = "34"
$ws.Range("J11").Value = "dfa"
